# Update per-unit voltage magnitude results (vm_pu) for the 380 kV case:
# slack bus voltage setpoint changed from 1.05 to 1.02 p.u. (column B),
# which propagates new load-flow solution values into columns C-F and I-N
# for every data row (rows 2-25). Columns G/H are unaffected and left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.035626920226119
$ws.Cells.Item(2, 4).Value = 1.043624602210796
$ws.Cells.Item(2, 5).Value = 1.044091706489379
$ws.Cells.Item(2, 6).Value = 1.05366442745138
$ws.Cells.Item(2, 9).Value = 1.032575100587231
$ws.Cells.Item(2, 10).Value = 1.040739663707883
$ws.Cells.Item(2, 11).Value = 1.046397942509901
$ws.Cells.Item(2, 12).Value = 1.046863731298549
$ws.Cells.Item(2, 13).Value = 1.056409771367323
$ws.Cells.Item(2, 14).Value = 1.01752394591051
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.036491012469765
$ws.Cells.Item(3, 4).Value = 1.044395999858868
$ws.Cells.Item(3, 5).Value = 1.044864885175132
$ws.Cells.Item(3, 6).Value = 1.054531073377579
$ws.Cells.Item(3, 9).Value = 1.032646022446604
$ws.Cells.Item(3, 10).Value = 1.041247813820245
$ws.Cells.Item(3, 11).Value = 1.046980516733484
$ws.Cells.Item(3, 12).Value = 1.047448176880773
$ws.Cells.Item(3, 13).Value = 1.057089369788153
$ws.Cells.Item(3, 14).Value = 1.017693974299031
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037050763424595
$ws.Cells.Item(4, 4).Value = 1.044896071442355
$ws.Cells.Item(4, 5).Value = 1.045366162283228
$ws.Cells.Item(4, 6).Value = 1.055093043636203
$ws.Cells.Item(4, 9).Value = 1.032690519455625
$ws.Cells.Item(4, 10).Value = 1.041576592950654
$ws.Cells.Item(4, 11).Value = 1.047357745120445
$ws.Cells.Item(4, 12).Value = 1.047826669408872
$ws.Cells.Item(4, 13).Value = 1.057529663920271
$ws.Cells.Item(4, 14).Value = 1.017803940699223
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.037286230835973
$ws.Cells.Item(5, 4).Value = 1.045106521299001
$ws.Cells.Item(5, 5).Value = 1.045577131646819
$ws.Cells.Item(5, 6).Value = 1.055329579445941
$ws.Cells.Item(5, 9).Value = 1.032708891579245
$ws.Cells.Item(5, 10).Value = 1.041714803657415
$ws.Cells.Item(5, 11).Value = 1.047516393631427
$ws.Cells.Item(5, 12).Value = 1.047985862113665
$ws.Cells.Item(5, 13).Value = 1.057714893237845
$ws.Cells.Item(5, 14).Value = 1.017850157311674
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037325775488237
$ws.Cells.Item(6, 4).Value = 1.045141869588943
$ws.Cells.Item(6, 5).Value = 1.045612567909207
$ws.Cells.Item(6, 6).Value = 1.055369311417888
$ws.Cells.Item(6, 9).Value = 1.032711956710249
$ws.Cells.Item(6, 10).Value = 1.041738009332726
$ws.Cells.Item(6, 11).Value = 1.047543034989323
$ws.Cells.Item(6, 12).Value = 1.048012595591737
$ws.Cells.Item(6, 13).Value = 1.057746001610916
$ws.Cells.Item(6, 14).Value = 1.017857916492216
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037053909171745
$ws.Cells.Item(7, 4).Value = 1.044898882622458
$ws.Cells.Item(7, 5).Value = 1.045368980355091
$ws.Cells.Item(7, 6).Value = 1.055096203127707
$ws.Cells.Item(7, 9).Value = 1.032690766259691
$ws.Cells.Item(7, 10).Value = 1.041578439762161
$ws.Cells.Item(7, 11).Value = 1.047359864749051
$ws.Cells.Item(7, 12).Value = 1.047828796258926
$ws.Cells.Item(7, 13).Value = 1.05753213845606
$ws.Cells.Item(7, 14).Value = 1.017804558300592
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.035918813749614
$ws.Cells.Item(8, 4).Value = 1.043885106864295
$ws.Cells.Item(8, 5).Value = 1.044352802062566
$ws.Cells.Item(8, 6).Value = 1.053957066435467
$ws.Cells.Item(8, 9).Value = 1.032599357244109
$ws.Cells.Item(8, 10).Value = 1.040911400401632
$ws.Cells.Item(8, 11).Value = 1.046594770660754
$ws.Cells.Item(8, 12).Value = 1.047061180859563
$ws.Cells.Item(8, 13).Value = 1.056639330270852
$ws.Cells.Item(8, 14).Value = 1.017581418547747
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.033923493451611
$ws.Cells.Item(9, 4).Value = 1.042105876478909
$ws.Cells.Item(9, 5).Value = 1.042569744663843
$ws.Cells.Item(9, 6).Value = 1.051958980511702
$ws.Cells.Item(9, 9).Value = 1.032427633894362
$ws.Cells.Item(9, 10).Value = 1.039735836973148
$ws.Cells.Item(9, 11).Value = 1.045248663912828
$ws.Cells.Item(9, 12).Value = 1.045711038984046
$ws.Cells.Item(9, 13).Value = 1.05507036206632
$ws.Cells.Item(9, 14).Value = 1.017187831090617
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.032596646662624
$ws.Cells.Item(10, 4).Value = 1.04092465569426
$ws.Cells.Item(10, 5).Value = 1.04138624599982
$ws.Cells.Item(10, 6).Value = 1.050633237420432
$ws.Cells.Item(10, 9).Value = 1.03230603060926
$ws.Cells.Item(10, 10).Value = 1.038952103916937
$ws.Cells.Item(10, 11).Value = 1.044352750838056
$ws.Cells.Item(10, 12).Value = 1.044812709893414
$ws.Cells.Item(10, 13).Value = 1.054027355152233
$ws.Cells.Item(10, 14).Value = 1.016925208826512
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032022926309189
$ws.Cells.Item(11, 4).Value = 1.040414367435197
$ws.Cells.Item(11, 5).Value = 1.040875036285391
$ws.Cells.Item(11, 6).Value = 1.050060698203215
$ws.Cells.Item(11, 9).Value = 1.032251694167479
$ws.Cells.Item(11, 10).Value = 1.038612749992008
$ws.Cells.Item(11, 11).Value = 1.04396518297115
$ws.Cells.Item(11, 12).Value = 1.044424160858728
$ws.Cells.Item(11, 13).Value = 1.053576447289667
$ws.Cells.Item(11, 14).Value = 1.016811441976499
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03180994481074
$ws.Cells.Item(12, 4).Value = 1.040225004054361
$ws.Cells.Item(12, 5).Value = 1.040685340386378
$ws.Cells.Item(12, 6).Value = 1.049848261440063
$ws.Cells.Item(12, 9).Value = 1.032231259216615
$ws.Cells.Item(12, 10).Value = 1.038486701339081
$ws.Cells.Item(12, 11).Value = 1.043821279859903
$ws.Cells.Item(12, 12).Value = 1.044279903061213
$ws.Cells.Item(12, 13).Value = 1.053409070033019
$ws.Cells.Item(12, 14).Value = 1.016769176936979
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.031855624437335
$ws.Cells.Item(13, 4).Value = 1.040265614948504
$ws.Cells.Item(13, 5).Value = 1.040726022165471
$ws.Cells.Item(13, 6).Value = 1.049893819416283
$ws.Cells.Item(13, 9).Value = 1.03223565397521
$ws.Cells.Item(13, 10).Value = 1.038513739067952
$ws.Cells.Item(13, 11).Value = 1.043852144969021
$ws.Cells.Item(13, 12).Value = 1.044310843809935
$ws.Cells.Item(13, 13).Value = 1.053444968018103
$ws.Cells.Item(13, 14).Value = 1.016778243238672
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032005318657443
$ws.Cells.Item(14, 4).Value = 1.040398710908021
$ws.Cells.Item(14, 5).Value = 1.040859352073213
$ws.Cells.Item(14, 6).Value = 1.050043133414866
$ws.Cells.Item(14, 9).Value = 1.032250010144637
$ws.Cells.Item(14, 10).Value = 1.038602330712117
$ws.Cells.Item(14, 11).Value = 1.043953286726657
$ws.Cells.Item(14, 12).Value = 1.044412235096591
$ws.Cells.Item(14, 13).Value = 1.053562609579758
$ws.Cells.Item(14, 14).Value = 1.016807948474376
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032097566658678
$ws.Cells.Item(15, 4).Value = 1.040480739682702
$ws.Cells.Item(15, 5).Value = 1.040941526283378
$ws.Cells.Item(15, 6).Value = 1.050135161206851
$ws.Cells.Item(15, 9).Value = 1.032258822084082
$ws.Cells.Item(15, 10).Value = 1.038656915320338
$ws.Cells.Item(15, 11).Value = 1.044015611079858
$ws.Cells.Item(15, 12).Value = 1.044474714486572
$ws.Cells.Item(15, 13).Value = 1.053635107054766
$ws.Cells.Item(15, 14).Value = 1.016826249942302
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.032634739862522
$ws.Cells.Item(16, 4).Value = 1.040958547039882
$ws.Cells.Item(16, 5).Value = 1.041420199866691
$ws.Cells.Item(16, 6).Value = 1.050671267078948
$ws.Cells.Item(16, 9).Value = 1.032309601368781
$ws.Cells.Item(16, 10).Value = 1.038974625996247
$ws.Cells.Item(16, 11).Value = 1.044378480331092
$ws.Cells.Item(16, 12).Value = 1.044838505872457
$ws.Cells.Item(16, 13).Value = 1.054057295756054
$ws.Cells.Item(16, 14).Value = 1.016932758153505
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.032971913252483
$ws.Cells.Item(17, 4).Value = 1.041258582407555
$ws.Cells.Item(17, 5).Value = 1.041720795968347
$ws.Cells.Item(17, 6).Value = 1.051007959486965
$ws.Cells.Item(17, 9).Value = 1.03234100415896
$ws.Cells.Item(17, 10).Value = 1.039173920516027
$ws.Cells.Item(17, 11).Value = 1.044606198318987
$ws.Cells.Item(17, 12).Value = 1.045066819688999
$ws.Cells.Item(17, 13).Value = 1.054322317793547
$ws.Cells.Item(17, 14).Value = 1.016999554991937
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.033168659125375
$ws.Cells.Item(18, 4).Value = 1.041433702472153
$ws.Cells.Item(18, 5).Value = 1.041896249359466
$ws.Cells.Item(18, 6).Value = 1.0512044926388
$ws.Cells.Item(18, 9).Value = 1.032359158600115
$ws.Cells.Item(18, 10).Value = 1.039290166285956
$ws.Cells.Item(18, 11).Value = 1.044739057813835
$ws.Cells.Item(18, 12).Value = 1.045200032990529
$ws.Cells.Item(18, 13).Value = 1.054476970119228
$ws.Cells.Item(18, 14).Value = 1.017038511631279
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.033235757621712
$ws.Cells.Item(19, 4).Value = 1.041493433254016
$ws.Cells.Item(19, 5).Value = 1.04195609485873
$ws.Cells.Item(19, 6).Value = 1.05127153009658
$ws.Cells.Item(19, 9).Value = 1.032365321251358
$ws.Cells.Item(19, 10).Value = 1.039329803158225
$ws.Cells.Item(19, 11).Value = 1.044784365434596
$ws.Cells.Item(19, 12).Value = 1.045245462313705
$ws.Cells.Item(19, 13).Value = 1.054529714295494
$ws.Cells.Item(19, 14).Value = 1.017051794004784
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.032935729623955
$ws.Cells.Item(20, 4).Value = 1.041226379594039
$ws.Cells.Item(20, 5).Value = 1.041688532345868
$ws.Cells.Item(20, 6).Value = 1.050971820444955
$ws.Cells.Item(20, 9).Value = 1.032337651717042
$ws.Cells.Item(20, 10).Value = 1.039152538023326
$ws.Cells.Item(20, 11).Value = 1.044581762656882
$ws.Cells.Item(20, 12).Value = 1.045042319451592
$ws.Cells.Item(20, 13).Value = 1.054293876240911
$ws.Cells.Item(20, 14).Value = 1.016992388818959
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.031961234010434
$ws.Cells.Item(21, 4).Value = 1.040359512470159
$ws.Cells.Item(21, 5).Value = 1.040820084474215
$ws.Cells.Item(21, 6).Value = 1.049999157804414
$ws.Cells.Item(21, 9).Value = 1.032245789561244
$ws.Cells.Item(21, 10).Value = 1.038576242600991
$ws.Cells.Item(21, 11).Value = 1.043923501418499
$ws.Cells.Item(21, 12).Value = 1.044382376038845
$ws.Cells.Item(21, 13).Value = 1.053527964021988
$ws.Cells.Item(21, 14).Value = 1.016799201211147
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.031349246937011
$ws.Cells.Item(22, 4).Value = 1.039815523374501
$ws.Cells.Item(22, 5).Value = 1.040275157919899
$ws.Cells.Item(22, 6).Value = 1.049388936701265
$ws.Cells.Item(22, 9).Value = 1.032186574631713
$ws.Cells.Item(22, 10).Value = 1.038213917956103
$ws.Cells.Item(22, 11).Value = 1.043509956116915
$ws.Cells.Item(22, 12).Value = 1.043967829573978
$ws.Cells.Item(22, 13).Value = 1.053047041867309
$ws.Cells.Item(22, 14).Value = 1.016677696437486
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.031673604354912
$ws.Cells.Item(23, 4).Value = 1.040103802631546
$ws.Cells.Item(23, 5).Value = 1.040563928792324
$ws.Cells.Item(23, 6).Value = 1.04971229967343
$ws.Cells.Item(23, 9).Value = 1.032218103532089
$ws.Cells.Item(23, 10).Value = 1.038405991216846
$ws.Cells.Item(23, 11).Value = 1.04372915255598
$ws.Cells.Item(23, 12).Value = 1.044187551394336
$ws.Cells.Item(23, 13).Value = 1.053301926836846
$ws.Cells.Item(23, 14).Value = 1.016742112047945
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.032952079206313
$ws.Cells.Item(24, 4).Value = 1.041240930305734
$ws.Cells.Item(24, 5).Value = 1.041703110515175
$ws.Cells.Item(24, 6).Value = 1.050988149671072
$ws.Cells.Item(24, 9).Value = 1.032339167042719
$ws.Cells.Item(24, 10).Value = 1.03916219984954
$ws.Cells.Item(24, 11).Value = 1.044592803971432
$ws.Cells.Item(24, 12).Value = 1.045053389925569
$ws.Cells.Item(24, 13).Value = 1.054306727539898
$ws.Cells.Item(24, 14).Value = 1.016995626919038
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.034438745294806
$ws.Cells.Item(25, 4).Value = 1.042564988843725
$ws.Cells.Item(25, 5).Value = 1.043029798211897
$ws.Cells.Item(25, 6).Value = 1.052474429204321
$ws.Cells.Item(25, 9).Value = 1.032473286466476
$ws.Cells.Item(25, 10).Value = 1.040039758402683
$ws.Cells.Item(25, 11).Value = 1.045596408368547
$ws.Cells.Item(25, 12).Value = 1.046059778304667
$ws.Cells.Item(25, 13).Value = 1.055475461339676
$ws.Cells.Item(25, 14).Value = 1.017289625680567